$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 422, shifting existing rows 422:443 down to 423:444.
$ws.Rows.Item(422).Insert()

# Populate the newly inserted row 422 with the new record's data.
$ws.Range("A422").Value = 10
$ws.Range("B422").Value = "Vega Modelo de Temuco"
$ws.Range("C422").Value = "La Araucanía"
$ws.Range("D422").Value = "3/23/2023"
$ws.Range("E422").Value = 9
$ws.Range("F422").Value = 100114013
$ws.Range("G422").Value = "Zanahoria"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 185
$ws.Range("K422").Value = 6000
$ws.Range("L422").Value = 7000
$ws.Range("M422").Value = 6324
$ws.Range("N422").Value = "$/saco 25 kilos"
$ws.Range("O422").Value = "Región de La Araucanía"
$ws.Range("P422").Value = 253
$ws.Range("Q422").Value = 25
$ws.Range("R422").Value = "Hortaliza"
